# The commit inserts one new data row (a weekly price observation) into the
# "Fruta, Feria Lagunitas de Puerto Montt - Frutilla" sheet, right before
# what used to be row 107. Every row from the old row 107 onward therefore
# shifts down by one (old 107->108, ..., old 116->117), and the new row 107
# is populated with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 107, pushing existing rows 107-116 down
# to 108-117.
$ws.Rows(107).Insert()

# Populate the newly inserted row 107 with the new observation.
$ws.Cells.Item(107, 1).Value  = 4
$ws.Cells.Item(107, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(107, 3).Value  = "Los Lagos"
$ws.Cells.Item(107, 4).Value  = 44505
$ws.Cells.Item(107, 5).Value  = 10
$ws.Cells.Item(107, 6).Value  = "Fruta"
$ws.Cells.Item(107, 7).Value  = 100101
$ws.Cells.Item(107, 8).Value  = "Berries"
$ws.Cells.Item(107, 9).Value  = 100112025
$ws.Cells.Item(107, 10).Value = "Frutilla"
$ws.Cells.Item(107, 11).Value = "Sin especificar"
$ws.Cells.Item(107, 12).Value = "Primera"
$ws.Cells.Item(107, 13).Value = 500
$ws.Cells.Item(107, 14).Value = 9000
$ws.Cells.Item(107, 15).Value = 10000
$ws.Cells.Item(107, 16).Value = 9500
$ws.Cells.Item(107, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(107, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(107, 19).Value = 1357
$ws.Cells.Item(107, 20).Value = 7
